$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 571
$ws.Cells.Item(3, 6).Value = 266
$ws.Cells.Item(5, 6).Value = 1425
$ws.Cells.Item(6, 6).Value = 734
$ws.Cells.Item(7, 6).Value = 371
$ws.Cells.Item(8, 6).Value = 51
$ws.Cells.Item(10, 6).Value = 6658
$ws.Cells.Item(13, 6).Value = 1923
$ws.Cells.Item(14, 6).Value = 4904
$ws.Cells.Item(16, 6).Value = 6075
$ws.Cells.Item(17, 6).Value = 7818
$ws.Cells.Item(19, 6).Value = 1102
$ws.Cells.Item(20, 6).Value = 788
$ws.Cells.Item(21, 6).Value = 4130
$ws.Cells.Item(22, 6).Value = 603
$ws.Cells.Item(23, 6).Value = 66
$ws.Cells.Item(28, 6).Value = 24
$ws.Cells.Item(29, 6).Value = 1541
$ws.Cells.Item(30, 6).Value = 608
$ws.Cells.Item(31, 6).Value = 745
$ws.Cells.Item(32, 6).Value = 1739
$ws.Cells.Item(34, 6).Value = 2017
$ws.Cells.Item(35, 6).Value = 251
$ws.Cells.Item(37, 6).Value = 1303
$ws.Cells.Item(39, 6).Value = 724
$ws.Cells.Item(40, 6).Value = 342
$ws.Cells.Item(41, 6).Value = 3804
$ws.Cells.Item(44, 6).Value = 373
$ws.Cells.Item(49, 6).Value = 3991

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 1311
$ws.Cells.Item(11, 6).Value = 3

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 4676

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 4676
$ws.Cells.Item(4, 6).Value = 571
$ws.Cells.Item(5, 6).Value = 1311
$ws.Cells.Item(8, 6).Value = 266
$ws.Cells.Item(11, 6).Value = 1425
$ws.Cells.Item(12, 6).Value = 734
$ws.Cells.Item(13, 6).Value = 51
$ws.Cells.Item(15, 6).Value = 6658
$ws.Cells.Item(16, 6).Value = 3
$ws.Cells.Item(18, 6).Value = 4904
$ws.Cells.Item(19, 6).Value = 6075
$ws.Cells.Item(20, 6).Value = 6075
$ws.Cells.Item(21, 6).Value = 7818
$ws.Cells.Item(23, 6).Value = 1102
$ws.Cells.Item(24, 6).Value = 788
$ws.Cells.Item(25, 6).Value = 4130
$ws.Cells.Item(26, 6).Value = 603
$ws.Cells.Item(27, 6).Value = 66
$ws.Cells.Item(31, 6).Value = 1541
$ws.Cells.Item(32, 6).Value = 608
$ws.Cells.Item(33, 6).Value = 745
$ws.Cells.Item(34, 6).Value = 1739
$ws.Cells.Item(36, 6).Value = 2017
$ws.Cells.Item(41, 6).Value = 724
$ws.Cells.Item(43, 6).Value = 342
$ws.Cells.Item(45, 6).Value = 3804
$ws.Cells.Item(47, 6).Value = 373
$ws.Cells.Item(51, 6).Value = 3991
